$wb = $excel.ActiveWorkbook

# --- ALC!row141 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 6448.875
$ws.Range("I141").Value = 2718.2
$ws.Range("K141").Value = 8154.599999999999
$ws.Range("M141").Value = -2974.599999999999

# --- ARM!row2 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 873.4
$ws.Range("I2").Value = 900
$ws.Range("J2").Value = 855.6667
$ws.Range("K2").Value = 900
$ws.Range("L2").Value = 855.6667
$ws.Range("M2").Value = -787
$ws.Range("N2").Value = -1081.6667

# --- ARM!row59 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H59").Value = 44759
$ws.Range("J59").Value = 44759
$ws.Range("L59").Value = 44759
$ws.Range("N59").Value = -46367

# --- ARM!row76 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 39750
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 39750
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 39750
$ws.Range("M76").Value = ""
$ws.Range("N76").Value = -40426

# --- ARM!row79 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 39750
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 39750
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 39750
$ws.Range("M79").Value = ""
$ws.Range("N79").Value = -42090

# --- ARM!row116 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 873.4
$ws.Range("I116").Value = 900
$ws.Range("J116").Value = 855.6667
$ws.Range("K116").Value = 900
$ws.Range("L116").Value = 855.6667
$ws.Range("M116").Value = 1394
$ws.Range("N116").Value = -5443.6667

# --- ARM!row122 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1343.9
$ws.Range("I122").Value = 1162.4166
$ws.Range("J122").Value = 1616.125
$ws.Range("K122").Value = 3487.2498
$ws.Range("L122").Value = 4848.375
$ws.Range("M122").Value = -1037.2498
$ws.Range("N122").Value = -9748.375

# --- BSM!row3 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 873.4
$ws.Range("I3").Value = 900
$ws.Range("J3").Value = 855.6667
$ws.Range("K3").Value = 900
$ws.Range("L3").Value = 855.6667
$ws.Range("M3").Value = -786
$ws.Range("N3").Value = -1083.6667

# --- BSM!row55 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 24333.334
$ws.Range("J55").Value = 24333.334
$ws.Range("L55").Value = 24333.334
$ws.Range("N55").Value = -24879.334

# --- CRP!row88 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").Value = ""

# --- CRP!row91 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").Value = ""

# --- CRP!row107 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 537.0244
$ws.Range("I107").Value = 367.4516
$ws.Range("J107").Value = 1062.7
$ws.Range("K107").Value = 367.4516
$ws.Range("L107").Value = 1062.7
$ws.Range("M107").Value = 1552.5484
$ws.Range("N107").Value = -4902.7

# --- CRP!row132 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 8131491.5
$ws.Range("I132").Value = 994.3
$ws.Range("K132").Value = 2982.9
$ws.Range("M132").Value = -452.8999999999996

# --- CUL!row106 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 1490
$ws.Range("J106").Value = 2000
$ws.Range("L106").Value = 6000
$ws.Range("N106").Value = -7892

# --- GSM!row29 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 1980
$ws.Range("J29").Value = 1980
$ws.Range("L29").Value = 1980
$ws.Range("N29").Value = -2560

# --- GSM!row80 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 14287943
$ws.Range("I80").Value = 2520
$ws.Range("K80").Value = 2520
$ws.Range("M80").Value = -1522

# --- GSM!row83 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 14287943
$ws.Range("I83").Value = 2520
$ws.Range("K83").Value = 12600
$ws.Range("M83").Value = -7608

# --- LTW!row42 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 3500
$ws.Range("J42").Value = 3500
$ws.Range("L42").Value = 3500
$ws.Range("N42").Value = -4626

# --- LTW!row49 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H49").Value = 3500
$ws.Range("J49").Value = 3500
$ws.Range("L49").Value = 3500
$ws.Range("N49").Value = -3794

# --- LTW!row61 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2047.8334
$ws.Range("I61").Value = 1997.4
$ws.Range("J61").Value = 2300
$ws.Range("K61").Value = 1997.4
$ws.Range("L61").Value = 2300
$ws.Range("M61").Value = -1795.4
$ws.Range("N61").Value = -2704

# --- LTW!row64 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 18000
$ws.Range("J64").Value = 18000
$ws.Range("L64").Value = 18000
$ws.Range("N64").Value = -18450

# --- LTW!row67 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H67").Value = 18000
$ws.Range("J67").Value = 18000
$ws.Range("L67").Value = 18000
$ws.Range("N67").Value = -19560

# --- LTW!row113 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2047.8334
$ws.Range("I113").Value = 1997.4
$ws.Range("J113").Value = 2300
$ws.Range("K113").Value = 1997.4
$ws.Range("L113").Value = 2300
$ws.Range("M113").Value = 172.5999999999999
$ws.Range("N113").Value = -6640

# --- LTW!row132 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 34640510
$ws.Range("I132").Value = 57144340
$ws.Range("J132").Value = 19234.846
$ws.Range("K132").Value = 171433020
$ws.Range("L132").Value = 57704.538
$ws.Range("M132").Value = -171430490
$ws.Range("N132").Value = -62764.538

# --- LTW!row136 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 153064930
$ws.Range("I136").Value = 114290590
$ws.Range("J136").Value = 250000750
$ws.Range("K136").Value = 342871770
$ws.Range("L136").Value = 750002250
$ws.Range("M136").Value = -342869220
$ws.Range("N136").Value = -750007350

# --- WVR!row16 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 59800
$ws.Range("J16").Value = 59800
$ws.Range("L16").Value = 59800
$ws.Range("N16").Value = -60384

# --- WVR!row63 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""

# --- WVR!row66 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""

# --- WVR!row82 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 29900
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 29900
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 29900
$ws.Range("M82").Value = ""
$ws.Range("N82").Value = -30666

# --- WVR!row85 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H85").Value = 29900
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 29900
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 29900
$ws.Range("M85").Value = ""
$ws.Range("N85").Value = -32552

# --- WVR!row100 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 833.3333
$ws.Range("I100").Value = 625
$ws.Range("K100").Value = 1250
$ws.Range("M100").Value = -709

# --- WVR!row120 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").Value = ""

# --- WVR!row121 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = ""

# --- WVR!row132 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 41946.25
$ws.Range("I132").Value = 104059.7
$ws.Range("J132").Value = 7438.778
$ws.Range("K132").Value = 312179.1
$ws.Range("L132").Value = 22316.334
$ws.Range("M132").Value = -309649.1
$ws.Range("N132").Value = -27376.334

# --- WVR!row136 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 19233372
$ws.Range("I136").Value = 45456430
$ws.Range("J136").Value = 3125.6667
$ws.Range("K136").Value = 136369290
$ws.Range("L136").Value = 9377.000100000001
$ws.Range("M136").Value = -136366740
$ws.Range("N136").Value = -14477.0001
